$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.9758064516129
$ws.Range("C2").Value = 2.07397622192867
$ws.Range("D2").Value = 2.03683241252302
$ws.Range("E2").Value = 1.9860248447205
$ws.Range("F2").Value = 1.85319712447976
